$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.02, 1.019094646265327, 1.025148302528871, 1.020272647394891, 1.032619258696602, 1.029025155640607, 1.02429935787862, 1.027974722202116, 1.023113422815218, 1.035423966993559, 1.025753980455587),
    @(1.02, 1.020036125399992, 1.025853833890716, 1.021071199257753, 1.034759587619546, 1.029229049537346, 1.024877485143739, 1.02848777156243, 1.023718199423944, 1.037369549619592, 1.026332928727735),
    @(1.02, 1.020645589142858, 1.026310265215647, 1.021588526313318, 1.036138934122311, 1.029359252625616, 1.025251251916586, 1.028818983791318, 1.024109490924874, 1.038622547520839, 1.02670722629227),
    @(1.02, 1.020901870865851, 1.026502126145679, 1.021806156251865, 1.036717504187739, 1.029413577135423, 1.025408306746106, 1.028958042465486, 1.024273980197612, 1.039147918337746, 1.026864504157662),
    @(1.02, 1.020944905399075, 1.02653433909859, 1.02184270580099, 1.036814572882794, 1.029422674285323, 1.025434672442061, 1.028981380306076, 1.024301598095739, 1.039236049756607, 1.026890907295929),
    @(1.02, 1.020649013343593, 1.026312828960596, 1.021591433723628, 1.036146670102674, 1.029359980133133, 1.025253350792476, 1.028820842618782, 1.024111688876931, 1.03862957297627, 1.026709328148804),
    @(1.02, 1.019412768940691, 1.025386759674863, 1.020542394989413, 1.033343768936809, 1.029094421498023, 1.024494804968849, 1.028148268162887, 1.02331781792056, 1.036082729809031, 1.025949705103101),
    @(1.02, 1.017236366766826, 1.023754187011795, 1.018698553742043, 1.028360489931845, 1.028613165965093, 1.023155699192395, 1.026957236602726, 1.021918623639666, 1.031548191839349, 1.024608697642847),
    @(1.02, 1.015786782401105, 1.022665323041222, 1.017472513039292, 1.025006604673706, 1.028283297828366, 1.022261312697127, 1.026159251316915, 1.020985640993043, 1.028491946729126, 1.023713041016169),
    @(1.02, 1.015159413943473, 1.022193717673224, 1.016942384589977, 1.023546372125335, 1.028138299762637, 1.021873640957981, 1.025812769532826, 1.020581606171808, 1.027160265774158, 1.023324818738687),
    @(1.02, 1.014926427750358, 1.022018524198061, 1.016745585050963, 1.023002739510728, 1.028084114302721, 1.021729582629283, 1.025683927664414, 1.020431522562637, 1.02666433623751, 1.023180555830637),
    @(1.02, 1.014976401986606, 1.022056104622318, 1.016787794067146, 1.023119407166117, 1.028095752081325, 1.021760486338843, 1.025711571167516, 1.020463716327162, 1.026770773491289, 1.02321150342701),
    @(1.02, 1.015140154287968, 1.022179236499359, 1.016926114756133, 1.023501460748712, 1.028133827450967, 1.021861734272491, 1.025802122346394, 1.020569200358492, 1.027119298486151, 1.023312895144337),
    @(1.02, 1.015241053620799, 1.02225509968866, 1.017011353782176, 1.023736691458786, 1.028157243595697, 1.021924108537253, 1.025857894936379, 1.020634191612214, 1.027333865017676, 1.023375357987709),
    @(1.02, 1.015828425304555, 1.02269661938976, 1.017507711857284, 1.025103344129028, 1.028292875151663, 1.022287032822232, 1.026182226133795, 1.02101245446002, 1.02858014794185, 1.023738797666801),
    @(1.02, 1.01619695109616, 1.022973540750658, 1.017819266588866, 1.025958446231257, 1.028377372921807, 1.022514579341268, 1.026385416103332, 1.021249716126444, 1.029359659179141, 1.023966667327994),
    @(1.02, 1.016411936002462, 1.023135052725545, 1.018001064046954, 1.026456445690358, 1.028426450465968, 1.022647265137724, 1.026503841860151, 1.021388102292243, 1.029813535883532, 1.0240995415535),
    @(1.02, 1.016485245406363, 1.023190122177898, 1.018063064629544, 1.026626121676486, 1.028443149314669, 1.022692501085388, 1.026544206476116, 1.021435287626085, 1.029968161663071, 1.024144841741399),
    @(1.02, 1.016157408661282, 1.022943830907654, 1.017785832176988, 1.025866781493667, 1.028368328697372, 1.022490169710077, 1.026363625212378, 1.021224260665783, 1.029276107904609, 1.023942223032327),
    @(1.02, 1.01509193203915, 1.022142977735935, 1.016885379617829, 1.023388989958041, 1.02812262423706, 1.021831920931736, 1.025775461236123, 1.020538138136107, 1.027016702296718, 1.023283039465219),
    @(1.02, 1.014422292778186, 1.021639343887389, 1.016319887732223, 1.021823926234665, 1.027966248876644, 1.021417707785826, 1.02540483093648, 1.020106704652715, 1.025588677584487, 1.022868238089131),
    @(1.02, 1.014777255924655, 1.021906339811625, 1.016619603178323, 1.022654289569007, 1.028049326276384, 1.021637322887156, 1.025601387770518, 1.020335419597578, 1.026346418536257, 1.023088165069094),
    @(1.02, 1.016175276093178, 1.022957255543007, 1.017800939523888, 1.025908203211841, 1.028372416037635, 1.022501199490196, 1.026373471861084, 1.021235762907147, 1.029313863596686, 1.023953268475999),
    @(1.02, 1.017798780324638, 1.024176331214086, 1.0191746711453, 1.02965422741026, 1.028739167568079, 1.023502180669101, 1.027265844123941, 1.022280382433029, 1.032726197203327, 1.024955671162993)
)

$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $rowData[$j]
    }
}
